# Populate the previously-empty "Survey 3" row (row 5) on Sheet1 with its
# Very Secure / Quite Secure / Neutral / Quite Insecure / Very Insecure counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 20

# Leave the selection where the author ended up after editing.
$ws.Range("C10").Select()
